# Implemented searching in products and in orders (tested partwise).
#
# Target sheet: "TODO CMS" (2nd tab) in the TODO workbook.
#  - Row 11 ("Order anhand orderID suchen können") flips from "offen" to "done".
#  - Three new TODO rows are appended (13-15) for the search feature work:
#      13: Bestellungen suchen   -> done
#      14: Produkte suchen       -> done
#      15: Doku schreiben        -> offen

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TODO CMS")

$xlPasteFormats = -4122

# --- Row 11: searching orders is now implemented -> done ---
$ws.Range("B11").Value = "done"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B11").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 13: Bestellungen suchen -> done ---
$ws.Range("A13").Value = "Bestellungen suchen"
$ws.Range("B13").Value = "done"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B13").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 14: Produkte suchen -> done ---
$ws.Range("A14").Value = "Produkte suchen"
$ws.Range("B14").Value = "done"
$ws.Range("B2").Copy() | Out-Null
$ws.Range("B14").PasteSpecial($xlPasteFormats) | Out-Null

# --- Row 15: Doku schreiben (still open) -> offen ---
$ws.Range("A15").Value = "Doku schreiben "
$ws.Range("B15").Value = "offen"
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B15").PasteSpecial($xlPasteFormats) | Out-Null

$excel.CutCopyMode = 0

$ws.Activate()
$ws.Range("B15").Select() | Out-Null
